# Actualización desde MV -datos-
# Applies the updated cell values to Sheet1 as captured by the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "I2"  = 6249
    "E5"  = 1340
    "I5"  = 10202
    "I9"  = 17075
    "G10" = 9141
    "I12" = 23247
    "G13" = 14067
    "I13" = 26444
    "I14" = 26915
    "C15" = 801
    "I15" = 29386
    "I16" = 31385
    "G17" = 16964
    "I17" = 33815
    "E18" = 3528
    "G18" = 17871
    "I18" = 34903
    "C19" = 941
    "E19" = 4082
    "G19" = 17228
    "I19" = 35825
    "C20" = 1060
    "E21" = 5331
    "G21" = 23820
    "I21" = 42133
    "I22" = 45495
    "G23" = 30554
    "I23" = 48871
    "E24" = 6421
    "G24" = 33344
    "I24" = 53866
    "G25" = 33465
    "I25" = 55360
    "I26" = 60158
    "I27" = 69618
    "C28" = 2315
    "E28" = 10877
    "G28" = 46331
    "I28" = 75163
    "G29" = 49924
    "I29" = 81365
    "C30" = 2655
    "E30" = 13041
    "G30" = 52205
    "I30" = 86551
    "E31" = 13846
    "I31" = 93016
    "E32" = 14087
    "G32" = 57804
    "I32" = 97792
    "C33" = 3006
    "E33" = 15152
    "I33" = 100405
    "C34" = 3032
    "E34" = 16264
    "G34" = 64044
    "I34" = 108822
    "C35" = 3322
    "E35" = 19027
    "G35" = 68619
    "I35" = 116841
    "C36" = 5102
    "E36" = 28390
    "G36" = 70212
    "I36" = 118723
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
